# Fix bug where the wrong positive phrase percentage was stored:
# append the missing trade row (row 16) that was dropped from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = 9569.9699999999993   # A16 Principle
$ws.Cells.Item($row, 2).Value = 9622.9                # B16 Start Principle
$ws.Cells.Item($row, 3).Value = 78.48                 # C16 BuyPrice
$ws.Cells.Item($row, 4).Value = 78.05                 # D16 SellPrice
$ws.Cells.Item($row, 5).Value = $false                # E16 IsShortSell
$ws.Cells.Item($row, 6).Value = -0.55000000000000004  # F16 Price Change %
$ws.Cells.Item($row, 7).Value = 42624.611145833333    # G16 Date (inherits date style from column)
$ws.Cells.Item($row, 8).Value = $false                # H16 Profitable
